$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Percentage of Outcomes" column with the cumulative binomial
# distribution formula, entered once in B5 and filled down through B25 so
# Excel records it as a single shared formula (t="shared").
$ws.Range("B5:B25").Formula = "=_xlfn.BINOM.DIST(A5, `$B`$2, `$B`$1,TRUE)"

# Leave the selection on B6, matching where the user clicked after filling
# the column down.
$ws.Range("B6").Select()
